$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = ""
$ws.Range("F2").Value = ""

# Row 3
$ws.Range("E3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F3").Value = ""

# Row 4
$ws.Range("E4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F4").Value = ""

# Row 5
$ws.Range("E5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F5").Value = ""

# Row 6
$ws.Range("E6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 7
$ws.Range("F7").Value = ""

# Row 8
$ws.Range("E8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F8").Value = ""

# Row 9
$ws.Range("E9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F9").Value = ""

# Row 10
$ws.Range("E10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F10").Value = ""

# Row 11
$ws.Range("E11").Value = "3,04 TL - 6,09 TL - 76,17 TL"

# Row 12
$ws.Range("C12").Value = "WU: 1.000,01 USD–"

# Row 13
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
$ws.Range("F13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"

# Row 14
$ws.Range("E14").Value = "1.660 TL - 1.660 TL"
$ws.Range("F14").Value = ""
